$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @("ERKEK YAKALI NUBUK CEKET SİYAH", "440 TL", "Ceket", "YAKALINUBUKSİYAH.jpg", "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır.", "Var"),
    @("ERKEK YAKALI NUBUK CEKET GRİ", "440 TL", "Ceket", "YAKALINUBUKGRİ.jpg", "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır.", "Var"),
    @("ERKEK YAKALI NUBUK CEKET LACİVERT", "440 TL", "Ceket", "YAKALINUBUKLACİ.jpg", "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır.", "Var")
)

$startRow = 95
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    # Column D (gorsel) is entered before column A (urun_adi) to match
    # the original authoring order of the shared-string table.
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

$ws.Range("I92").Select()
